$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values for rows 2-6 (columns D..AJ)
# Row 2
$ws.Range("D2").Value = 1574
$ws.Range("E2").Value = 82
$ws.Range("F2").Value = 82
$ws.Range("G2").Value = 54
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 1722
$ws.Range("L2").Value = 739
$ws.Range("M2").Value = 983
$ws.Range("N2").Value = 969
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 129
$ws.Range("R2").Value = -139
$ws.Range("S2").Value = 24
$ws.Range("T2").Value = 142
$ws.Range("U2").Value = -13
$ws.Range("V2").Value = 437
$ws.Range("W2").Value = 5.19
$ws.Range("X2").Value = 2.13
$ws.Range("Y2").Value = 3.84
$ws.Range("Z2").Value = 1.97
$ws.Range("AA2").Value = 75.19
$ws.Range("AB2").Value = 855.11
$ws.Range("AC2").Value = 184
$ws.Range("AD2").Value = 20.42
$ws.Range("AE2").Value = 4969
$ws.Range("AF2").Value = 0.76
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 2.66
$ws.Range("AI2").Value = 52.88
$ws.Range("AJ2").Value = 20000000

# Row 3
$ws.Range("D3").Value = 1574
$ws.Range("E3").Value = 138
$ws.Range("F3").Value = 138
$ws.Range("G3").Value = 131
$ws.Range("H3").Value = 89
$ws.Range("I3").Value = 91
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 1727
$ws.Range("L3").Value = 670
$ws.Range("M3").Value = 1058
$ws.Range("N3").Value = 1046
$ws.Range("O3").Value = 12
$ws.Range("P3").Value = 100
$ws.Range("Q3").Value = 235
$ws.Range("R3").Value = -67
$ws.Range("S3").Value = -104
$ws.Range("T3").Value = 76
$ws.Range("U3").Value = 159
$ws.Range("V3").Value = 352
$ws.Range("W3").Value = 8.76
$ws.Range("X3").Value = 5.66
$ws.Range("Y3").Value = 9.06
$ws.Range("Z3").Value = 5.16
$ws.Range("AA3").Value = 63.31
$ws.Range("AB3").Value = 929.95
$ws.Range("AC3").Value = 456
$ws.Range("AD3").Value = 10.56
$ws.Range("AE3").Value = 5365
$ws.Range("AF3").Value = 0.9
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 2.07
$ws.Range("AI3").Value = 21.37
$ws.Range("AJ3").Value = 20000000

# Row 4
$ws.Range("D4").Value = 1578
$ws.Range("E4").Value = 163
$ws.Range("F4").Value = 163
$ws.Range("G4").Value = 161
$ws.Range("H4").Value = 112
$ws.Range("I4").Value = 109
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1822
$ws.Range("L4").Value = 676
$ws.Range("M4").Value = 1146
$ws.Range("N4").Value = 1131
$ws.Range("O4").Value = 15
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 203
$ws.Range("R4").Value = -41
$ws.Range("S4").Value = -84
$ws.Range("T4").Value = 46
$ws.Range("U4").Value = 157
$ws.Range("V4").Value = 288
$ws.Range("W4").Value = 10.33
$ws.Range("X4").Value = 7.11
$ws.Range("Y4").Value = 10.02
$ws.Range("Z4").Value = 6.32
$ws.Range("AA4").Value = 58.96
$ws.Range("AB4").Value = 1020.15
$ws.Range("AC4").Value = 545
$ws.Range("AD4").Value = 11.9
$ws.Range("AE4").Value = 5802
$ws.Range("AF4").Value = 1.12
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.54
$ws.Range("AI4").Value = 17.88
$ws.Range("AJ4").Value = 20000000

# Row 5
$ws.Range("D5").Value = 1646
$ws.Range("E5").Value = 125
$ws.Range("F5").Value = 125
$ws.Range("G5").Value = 380
$ws.Range("H5").Value = 296
$ws.Range("I5").Value = 293
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1792
$ws.Range("L5").Value = 377
$ws.Range("M5").Value = 1415
$ws.Range("N5").Value = 1397
$ws.Range("O5").Value = 17
$ws.Range("P5").Value = 100
$ws.Range("Q5").Value = 151
$ws.Range("R5").Value = 232
$ws.Range("S5").Value = -261
$ws.Range("T5").Value = 52
$ws.Range("U5").Value = 99
$ws.Range("V5").Value = 46
$ws.Range("W5").Value = 7.58
$ws.Range("X5").Value = 17.96
$ws.Range("Y5").Value = 23.16
$ws.Range("Z5").Value = 16.37
$ws.Range("AA5").Value = 26.67
$ws.Range("AB5").Value = 1296.47
$ws.Range("AC5").Value = 1464
$ws.Range("AD5").Value = 3.33
$ws.Range("AE5").Value = 7165
$ws.Range("AF5").Value = 0.68
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 2.05
$ws.Range("AI5").Value = 6.66
$ws.Range("AJ5").Value = 20000000

# Row 6
$ws.Range("D6").Value = 1439
$ws.Range("E6").Value = 72
$ws.Range("F6").Value = 72
$ws.Range("G6").Value = 81
$ws.Range("H6").Value = 57
$ws.Range("I6").Value = 57
$ws.Range("K6").Value = 1761
$ws.Range("L6").Value = 326
$ws.Range("M6").Value = 1435
$ws.Range("N6").Value = 1419
$ws.Range("P6").Value = 100
$ws.Range("Q6").Value = 141
$ws.Range("R6").Value = -187
$ws.Range("S6").Value = -21
$ws.Range("T6").Value = 195
$ws.Range("U6").Value = -54
$ws.Range("V6").Value = 46
$ws.Range("W6").Value = 5.03
$ws.Range("X6").Value = 3.93
$ws.Range("Y6").Value = 4.07
$ws.Range("Z6").Value = 3.19
$ws.Range("AA6").Value = 22.69
$ws.Range("AB6").Value = 1324.63
$ws.Range("AC6").Value = 287
$ws.Range("AD6").Value = 13.01
$ws.Range("AE6").Value = 7278
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 2.68
$ws.Range("AI6").Value = 34.02
$ws.Range("AJ6").Value = 20000000

# Clear cells D..AJ for rows 7-9 (data no longer present in source)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
